# ---------------------------------------------------------------------------
# Update latest output (run 213)
# Applies the "optimisation_result.xlsx" refresh: Schedule totals + a new
# completed pump-run row, and the Detailed time series rolling forward a day
# (older forecast rows become historical with actualised prices, and a new
# day of forecast rows is appended).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Schedule")
$ws2 = $wb.Worksheets.Item("Detailed")

$dtFormat = "YYYY-MM-DD HH:MM:SS"
$dateFormat = "YYYY-MM-DD"

# --- Schedule: recompute cost/unit-cost for the existing pump run (row 2) ---
$ws1.Cells.Item(2,5).Value = 358.3358902500001
$ws1.Cells.Item(2,6).Value = 8.617986778499281

# --- Schedule: append the newly-completed pump run (row 3) ---
$ws1.Cells.Item(3,1).Value = 46056.20833333334
$ws1.Cells.Item(3,1).NumberFormat = $dtFormat
$ws1.Cells.Item(3,2).Value = 46056.66666666666
$ws1.Cells.Item(3,2).NumberFormat = $dtFormat
$ws1.Cells.Item(3,3).Value = 11
$ws1.Cells.Item(3,4).Value = 41.58
$ws1.Cells.Item(3,5).Value = 434.23648125
$ws1.Cells.Item(3,6).Value = 10.44339781746032

# --- Detailed: refresh Price / Type for rows that have now actualised (13-49) ---
# Each tuple is (row, newPrice-or-$null, newType-or-$null)
$rowUpdates = @(
  @(13, 72.37465, $null),
  @(14, 73.19, $null),
  @(15, 73.19, "historical"),
  @(16, 57.31, "historical"),
  @(17, $null, "historical"),
  @(18, $null, "historical"),
  @(19, 31.41464, "historical"),
  @(20, 22.07, "historical"),
  @(21, 0.7, "historical"),
  @(22, 0.02795, "historical"),
  @(23, 0.02799, "historical"),
  @(24, $null, "historical"),
  @(25, 0.0085, "historical"),
  @(26, -5.01, "historical"),
  @(27, -12.01, "historical"),
  @(28, -12.01, "historical"),
  @(29, -14, "historical"),
  @(30, $null, "historical"),
  @(31, $null, "historical"),
  @(32, -8.94981, "historical"),
  @(33, -7.83618, "historical"),
  @(34, -9.5, $null),
  @(35, -7.60218, $null),
  @(36, -7.82404, $null),
  @(37, -2.9022, $null),
  @(38, -2.26986, $null),
  @(39, 23.73987, $null),
  @(40, 73.19, $null),
  @(41, 84.79, $null),
  @(42, 92.65842, $null),
  @(43, 103.08647, $null),
  @(44, 108.89, $null),
  @(45, 105.79, $null),
  @(46, 84.79, $null),
  @(47, 88.05589, $null),
  @(48, 93.15783, $null),
  @(49, 97.75772, $null)
)

foreach ($u in $rowUpdates) {
  $r = $u[0]
  $bVal = $u[1]
  $cVal = $u[2]
  if ($bVal -ne $null) { $ws2.Cells.Item($r,2).Value = $bVal }
  if ($cVal -ne $null) { $ws2.Cells.Item($r,3).Value = $cVal }
}

# --- Detailed: append the next day of forecast rows (50-97) ---
# Each tuple is (row, DateTime, Price, Type, Date, Pump_Status)
$newRows = @(
  @(50, 46056, 84.79, "forecast", 46056, "OFF"),
  @(51, 46056.02083333334, 84.79, "forecast", 46056, "OFF"),
  @(52, 46056.04166666666, 67.08355, "forecast", 46056, "OFF"),
  @(53, 46056.0625, 61.04389, "forecast", 46056, "OFF"),
  @(54, 46056.08333333334, 57.06007, "forecast", 46056, "OFF"),
  @(55, 46056.10416666666, 57.06007, "forecast", 46056, "OFF"),
  @(56, 46056.125, 57.41437, "forecast", 46056, "OFF"),
  @(57, 46056.14583333334, 58.12899, "forecast", 46056, "OFF"),
  @(58, 46056.16666666666, 59.54902, "forecast", 46056, "OFF"),
  @(59, 46056.1875, 59.30428, "forecast", 46056, "OFF"),
  @(60, 46056.20833333334, 69.3127, "forecast", 46056, "ON"),
  @(61, 46056.22916666666, 84.79, "forecast", 46056, "ON"),
  @(62, 46056.25, 101.25, "forecast", 46056, "ON"),
  @(63, 46056.27083333334, 89.71073, "forecast", 46056, "ON"),
  @(64, 46056.29166666666, 63.62195, "forecast", 46056, "ON"),
  @(65, 46056.3125, 21.27537, "forecast", 46056, "ON"),
  @(66, 46056.33333333334, 9.91834, "forecast", 46056, "ON"),
  @(67, 46056.35416666666, 17.89224, "forecast", 46056, "ON"),
  @(68, 46056.375, 0.51, "forecast", 46056, "ON"),
  @(69, 46056.39583333334, 0.00921, "forecast", 46056, "ON"),
  @(70, 46056.41666666666, 0.00985, "forecast", 46056, "ON"),
  @(71, 46056.4375, 0, "forecast", 46056, "ON"),
  @(72, 46056.45833333334, -5.01, "forecast", 46056, "ON"),
  @(73, 46056.47916666666, -5.50985, "forecast", 46056, "ON"),
  @(74, 46056.5, -5.06903, "forecast", 46056, "ON"),
  @(75, 46056.52083333334, -5.50985, "forecast", 46056, "ON"),
  @(76, 46056.54166666666, -5.17224, "forecast", 46056, "ON"),
  @(77, 46056.5625, -5.17224, "forecast", 46056, "ON"),
  @(78, 46056.58333333334, -4.76643, "forecast", 46056, "ON"),
  @(79, 46056.60416666666, 0.51, "forecast", 46056, "ON"),
  @(80, 46056.625, 0.7, "forecast", 46056, "ON"),
  @(81, 46056.64583333334, 22.07, "forecast", 46056, "ON"),
  @(82, 46056.66666666666, 30.91077, "forecast", 46056, "OFF"),
  @(83, 46056.6875, 22.07, "forecast", 46056, "OFF"),
  @(84, 46056.70833333334, 46.6152, "forecast", 46056, "OFF"),
  @(85, 46056.72916666666, 58.27678, "forecast", 46056, "OFF"),
  @(86, 46056.75, 73.18616, "forecast", 46056, "OFF"),
  @(87, 46056.77083333334, 78, "forecast", 46056, "OFF"),
  @(88, 46056.79166666666, 84.79, "forecast", 46056, "OFF"),
  @(89, 46056.8125, 84.79, "forecast", 46056, "OFF"),
  @(90, 46056.83333333334, 83.66007, "forecast", 46056, "OFF"),
  @(91, 46056.85416666666, 84.79, "forecast", 46056, "OFF"),
  @(92, 46056.875, 84.79, "forecast", 46056, "OFF"),
  @(93, 46056.89583333334, 84.79, "forecast", 46056, "OFF"),
  @(94, 46056.91666666666, 79.11609, "forecast", 46056, "OFF"),
  @(95, 46056.9375, 84.79, "forecast", 46056, "OFF"),
  @(96, 46056.95833333334, 84.79, "forecast", 46056, "OFF"),
  @(97, 46056.97916666666, 78, "forecast", 46056, "OFF")
)

foreach ($nr in $newRows) {
  $r = $nr[0]
  $ws2.Cells.Item($r,1).Value = $nr[1]
  $ws2.Cells.Item($r,1).NumberFormat = $dtFormat
  $ws2.Cells.Item($r,2).Value = $nr[2]
  $ws2.Cells.Item($r,3).Value = $nr[3]
  $ws2.Cells.Item($r,4).Value = $nr[4]
  $ws2.Cells.Item($r,4).NumberFormat = $dateFormat
  $ws2.Cells.Item($r,5).Value = $nr[5]
}

Write-Host "Edit complete"